$d = $word.ActiveDocument

$replacements = @(
    @("440÷4=110, 0", "125÷7=17, 6"),
    @("379÷2=189, 1", "466÷3=155, 1"),
    @("817÷5=163, 2", "805÷2=402, 1"),
    @("299÷6=49, 5", "146÷9=16, 2"),
    @("983÷9=109, 2", "345÷3=115, 0"),
    @("441÷7=63, 0", "755÷7=107, 6"),
    @("471÷5=94, 1", "515÷3=171, 2"),
    @("216÷5=43, 1", "157÷3=52, 1"),
    @("878÷5=175, 3", "778÷9=86, 4"),
    @("968÷5=193, 3", "368÷3=122, 2"),
    @("555÷8=69, 3", "721÷8=90, 1"),
    @("621÷3=207, 0", "628÷7=89, 5"),
    @("123÷3=41, 0", "109÷2=54, 1"),
    @("894÷7=127, 5", "722÷4=180, 2"),
    @("546÷2=273, 0", "694÷3=231, 1"),
    @("738÷2=369, 0", "766÷4=191, 2"),
    @("960÷6=160, 0", "300÷3=100, 0"),
    @("557÷5=111, 2", "153÷3=51, 0"),
    @("711÷5=142, 1", "753÷3=251, 0"),
    @("688÷8=86, 0", "436÷9=48, 4"),
    @("317÷3=105, 2", "603÷4=150, 3"),
    @("321÷3=107, 0", "828÷7=118, 2"),
    @("967÷3=322, 1", "143÷2=71, 1"),
    @("618÷5=123, 3", "515÷7=73, 4"),
    @("455÷8=56, 7", "602÷3=200, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
